$d = $word.ActiveDocument

# 1. Update the date field "12/8/22" -> "2/17/23"
$d.Content.Find.Execute("12/8/22", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2/17/23", 2)

# 2. Update the quarto image descr paths from /usr/lib/quarto-cli to /opt/quarto
$d.Content.Find.Execute("/usr/lib/quarto-cli/share/formats/docx/", $true, $false, $false, $false, $false,
                         $true, 1, $false, "/opt/quarto/share/formats/docx/", 2)

# 3. Update "As of December 8, 2022" -> "As of February 17, 2023"
$d.Content.Find.Execute("As of December 8, 2022 there are roughly 6.3 pages of text.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "As of February 17, 2023 there are roughly 6.3 pages of text.", 2)

# 4. Insert a new paragraph with style "FirstParagraph" and text "P" after the
#    paragraph containing "5 1626 11344      total"
$found = $d.Content.Find.Execute("5 1626 11344      total", $true, $false, $false, $false, $false,
                                  $true, 1, $false, "", 0)

$para = $d.Paragraphs.Item($d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*5 1626 11344*total*") {
        $para = $p
        break
    }
}

$insertRange = $para.Range
$insertRange.Collapse(0)
$newPara = $insertRange.Paragraphs.Add($insertRange)
$newPara.Range.Text = "P"
$newPara.Range.Style = "FirstParagraph"
